# [FEAT] : zip file 만들기
#
# Paragraph 1: "성명 : dlfma"   -> "성명 : " + "+++= " + "name"    + "+++"
# Paragraph 2: "주소 : 주소"     -> "주소 : " + "+++= " + "address" + "+++"
#
# Each paragraph keeps its first run (the "label : " run) untouched and
# replaces everything after it (the old value run plus the two trailing
# empty runs) with three new runs: "+++= ", the new value, "+++".

$d = $word.ActiveDocument

function Set-FieldValue {
    param(
        [int]$ParaIndex,
        [string]$NewValue
    )

    $p = $d.Paragraphs($ParaIndex)

    # Locate the end of the "label : " run by searching for the literal
    # ": " separator inside this paragraph only.
    $labelRange = $p.Range.Duplicate
    $labelRange.Find.ClearFormatting()
    $found = $labelRange.Find.Execute(": ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Could not find label separator ': ' in paragraph $ParaIndex"
    }

    $valueStart = $labelRange.End
    $paraEnd = $p.Range.End - 1   # exclude the paragraph mark

    $replaceRange = $d.Range($valueStart, $paraEnd)

    $xml = '<w:p>' +
        '<w:r><w:rPr><w:shd w:val="clear" w:color="auto" w:fill="auto"/></w:rPr><w:t xml:space="preserve">+++= </w:t></w:r>' +
        '<w:r><w:rPr><w:shd w:val="clear" w:color="auto" w:fill="auto"/></w:rPr><w:t xml:space="preserve">' + $NewValue + '</w:t></w:r>' +
        '<w:r><w:rPr><w:shd w:val="clear" w:color="auto" w:fill="auto"/></w:rPr><w:t xml:space="preserve">+++</w:t></w:r>' +
        '</w:p>'

    $replaceRange.InsertXML($xml)
}

Set-FieldValue 1 "name"
Set-FieldValue 2 "address"
